# Update coefficients worksheet to reflect "moving averages" refit:
#  - F/G/H values (intercept/slope/quadratic coefficients) get new figures
#    for several rows
#  - The R^2 column (I) is dropped entirely since it is no longer computed
#  - Selection / view state updated

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("coefficients")

# --- Updated regression coefficients (column F, G, H) ---
# (literals written in plain-decimal form -- the PS parser here doesn't
#  accept scientific notation -- but they are the exact IEEE-754 doubles
#  from the target workbook)
$ws.Range("F3").Value = -2.3289429
$ws.Range("G3").Value = 0.0716877
$ws.Range("H3").Value = -0.0000914

$ws.Range("F4").Value = -2.3835485
$ws.Range("G4").Value = 0.0642864
$ws.Range("H4").Value = 0.0008434

$ws.Range("F5").Value = -2.418263
$ws.Range("G5").Value = 0.045934
$ws.Range("H5").Value = 0.003216

$ws.Range("F6").Value = -2.366372
$ws.Range("G6").Value = 0.008832
$ws.Range("H6").Value = 0.005009

$ws.Range("F10").Value = -2.134723
$ws.Range("G10").Value = 0.055158

$ws.Range("F11").Value = -2.116415
$ws.Range("G11").Value = 0.052711

$ws.Range("F12").Value = -2.06886
$ws.Range("G12").Value = 0.047427

$ws.Range("F13").Value = -2.115944
$ws.Range("G13").Value = 0.052842

# --- Drop the now-unused R^2 column entirely (header + values) ---
$ws.Columns.Item(9).Delete()

# Row 1 no longer needs its taller custom height once the sheet settles
# back down to a plain header row.
$ws.Rows.Item(1).AutoFit()

# --- View / selection bookkeeping to mirror the authored edit ---
$ws.Range("F20").Select()
